# Updated cryptos list on Fri May 31 09:50:26 UTC 2024 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns for the coinranking feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write a literal-text value into a cell without letting Excel's
# number/date auto-detection coerce numeric-looking strings (e.g. '1.00',
# '592.16') into real numbers -- the source feed stores these as plain text.
function Set-TextValue([object]$range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value2 = $text
    $range.Style = "Normal"
}

$updates = @(
    @{ Cell = "D2"; Text = "67.975.50" }
    @{ Cell = "E2"; Text = "  +0.52%  " }
    @{ Cell = "D3"; Text = "3.745.12" }
    @{ Cell = "E3"; Text = "  +0.49%  " }
    @{ Cell = "D4"; Text = "1.00" }
    @{ Cell = "E4"; Text = "  -0.04%  " }
    @{ Cell = "D5"; Text = "592.16" }
    @{ Cell = "E5"; Text = "  +0.19%  " }
    @{ Cell = "D6"; Text = "166.87" }
    @{ Cell = "E6"; Text = "  +0.92%  " }
    @{ Cell = "D7"; Text = "3.744.31" }
    @{ Cell = "E7"; Text = "  +0.49%  " }
    @{ Cell = "E8"; Text = "  -0.02%  " }
    @{ Cell = "D9"; Text = "0.520" }
    @{ Cell = "E9"; Text = "  +0.39%  " }
    @{ Cell = "E10"; Text = "  +0.37%  " }
    @{ Cell = "E11"; Text = "  -0.88%  " }
    @{ Cell = "E12"; Text = "  -0.43%  " }
    @{ Cell = "D13"; Text = "0.0000258" }
    @{ Cell = "E13"; Text = "  -1.04%  " }
    @{ Cell = "E14"; Text = "  -0.16%  " }
    @{ Cell = "D15"; Text = "4.368.49" }
    @{ Cell = "E15"; Text = "  +0.50%  " }
    @{ Cell = "D16"; Text = "3.706.85" }
    @{ Cell = "E16"; Text = "  -0.64%  " }
    @{ Cell = "D17"; Text = "67.930.17" }
    @{ Cell = "D18"; Text = "17.82" }
    @{ Cell = "E18"; Text = "  -2.21%  " }
    @{ Cell = "E19"; Text = "  -0.46%  " }
    @{ Cell = "E20"; Text = "  +0.49%  " }
    @{ Cell = "D21"; Text = "10.62" }
    @{ Cell = "E21"; Text = "  -0.43%  " }
    @{ Cell = "D22"; Text = "464.23" }
    @{ Cell = "E22"; Text = "  -0.52%  " }
    @{ Cell = "D23"; Text = "0.694" }
    @{ Cell = "E23"; Text = "  -0.55%  " }
    @{ Cell = "D24"; Text = "0.0000148" }
    @{ Cell = "E24"; Text = "  +10.60%  " }
    @{ Cell = "D25"; Text = "83.71" }
    @{ Cell = "E25"; Text = "  +1.22%  " }
    @{ Cell = "E26"; Text = "  +0.25%  " }
    @{ Cell = "D27"; Text = "11.81" }
    @{ Cell = "E27"; Text = "  -1.23%  " }
    @{ Cell = "E28"; Text = "  +0.22%  " }
    @{ Cell = "E29"; Text = "  +0.08%  " }
    @{ Cell = "E30"; Text = "  +0.02%  " }
    @{ Cell = "D31"; Text = "7.24" }
    @{ Cell = "E31"; Text = "  -1.02%  " }
    @{ Cell = "D32"; Text = "29.69" }
    @{ Cell = "E32"; Text = "  +0.68%  " }
    @{ Cell = "E33"; Text = "  -3.61%  " }
    @{ Cell = "D35"; Text = "9.11" }
    @{ Cell = "E35"; Text = "  +1.00%  " }
    @{ Cell = "D36"; Text = "3.698.77" }
    @{ Cell = "E36"; Text = "  +0.59%  " }
    @{ Cell = "E37"; Text = "  -0.49%  " }
    @{ Cell = "D38"; Text = "3.44" }
    @{ Cell = "E38"; Text = "  +0.71%  " }
    @{ Cell = "E39"; Text = "  -0.03%  " }
    @{ Cell = "D40"; Text = "0.995" }
    @{ Cell = "E40"; Text = "  +0.80%  " }
    @{ Cell = "E41"; Text = "  +0.24%  " }
    @{ Cell = "E42"; Text = "  -0.02%  " }
    @{ Cell = "D44"; Text = "44.12" }
    @{ Cell = "E44"; Text = "  +16.85%  " }
    @{ Cell = "E45"; Text = "  -1.67%  " }
    @{ Cell = "D46"; Text = "46.80" }
    @{ Cell = "E46"; Text = "  +3.55%  " }
    @{ Cell = "D47"; Text = "1.90" }
    @{ Cell = "E47"; Text = "  -0.14%  " }
    @{ Cell = "E48"; Text = "  -1.50%  " }
    @{ Cell = "D49"; Text = "145.08" }
    @{ Cell = "E49"; Text = "  +1.48%  " }
    @{ Cell = "D50"; Text = "388.35" }
    @{ Cell = "E50"; Text = "  +0.71%  " }
    @{ Cell = "D51"; Text = "2.764.48" }
    @{ Cell = "E51"; Text = "  +3.31%  " }
)

foreach ($u in $updates) {
    Set-TextValue $ws.Range($u.Cell) $u.Text
}
